$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as text, preserving exact string content (avoids numeric/date coercion)
function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2
Set-TextCell $ws 'D2' '64.385.58'
Set-TextCell $ws 'E2' '  -2.22%  '

# Row 3
Set-TextCell $ws 'D3' '2.615.00'
Set-TextCell $ws 'E3' '  -2.29%  '

# Row 4
Set-TextCell $ws 'E4' '  +0.00%  '

# Row 5
Set-TextCell $ws 'D5' '576.22'
Set-TextCell $ws 'E5' '  -4.04%  '

# Row 6
Set-TextCell $ws 'D6' '156.12'
Set-TextCell $ws 'E6' '  -1.04%  '

# Row 7
Set-TextCell $ws 'D7' '0.641'
Set-TextCell $ws 'E7' '  +4.98%  '

# Row 8
Set-TextCell $ws 'E8' '  +0.02%  '

# Row 9
Set-TextCell $ws 'E9' '  -6.29%  '

# Row 10
Set-TextCell $ws 'D10' '5.80'
Set-TextCell $ws 'E10' '  -0.87%  '

# Row 11
Set-TextCell $ws 'E11' '  -2.56%  '

# Row 12
Set-TextCell $ws 'D12' '0.155'
Set-TextCell $ws 'E12' '  +0.23%  '

# Row 13
Set-TextCell $ws 'D13' '28.33'
Set-TextCell $ws 'E13' '  -3.01%  '

# Row 14
Set-TextCell $ws 'E14' '  -8.38%  '

# Row 15
Set-TextCell $ws 'D15' '3.090.60'
Set-TextCell $ws 'E15' '  -2.15%  '

# Row 16
Set-TextCell $ws 'D16' '64.322.75'
Set-TextCell $ws 'E16' '  -2.12%  '

# Row 17
Set-TextCell $ws 'D17' '2.637.71'
Set-TextCell $ws 'E17' '  -1.61%  '

# Row 18
Set-TextCell $ws 'D18' '12.21'
Set-TextCell $ws 'E18' '  -4.02%  '

# Row 19
Set-TextCell $ws 'E19' '  -2.90%  '

# Row 20
Set-TextCell $ws 'D20' '345.90'
Set-TextCell $ws 'E20' '  -1.93%  '

# Row 21
Set-TextCell $ws 'E21' '  -3.04%  '

# Row 22
Set-TextCell $ws 'E22' '  -0.28%  '

# Row 23
Set-TextCell $ws 'D23' '68.15'
Set-TextCell $ws 'E23' '  -2.30%  '

# Row 24
Set-TextCell $ws 'B24' 'PEPE'
Set-TextCell $ws 'C24' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws 'D24' '0.0000109'
Set-TextCell $ws 'E24' '  -7.43%  '

# Row 25
Set-TextCell $ws 'B25' 'SuiNetwork'
Set-TextCell $ws 'C25' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell $ws 'D25' '1.71'
Set-TextCell $ws 'E25' '  +2.40%  '

# Row 26
Set-TextCell $ws 'D26' '9.31'
Set-TextCell $ws 'E26' '  -3.96%  '

# Row 27
Set-TextCell $ws 'D27' '1.58'
Set-TextCell $ws 'E27' '  -1.54%  '

# Row 28
Set-TextCell $ws 'B28' 'Bittensor'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D28' '546.10'
Set-TextCell $ws 'E28' '  +1.85%  '

# Row 29
Set-TextCell $ws 'B29' 'Kaspa'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D29' '0.162'
Set-TextCell $ws 'E29' '  -2.33%  '

# Row 30
Set-TextCell $ws 'B30' 'Binance-PegBSC-USD'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D30' '1.00'
Set-TextCell $ws 'E30' '  +0.41%  '

# Row 31
Set-TextCell $ws 'E31' '  -0.59%  '

# Row 32
Set-TextCell $ws 'D32' '2.08'
Set-TextCell $ws 'E32' '  -2.50%  '

# Row 33
Set-TextCell $ws 'D33' '1.73'
Set-TextCell $ws 'E33' '  -2.06%  '

# Row 34
Set-TextCell $ws 'D34' '6.37'
Set-TextCell $ws 'E34' '  -2.12%  '

# Row 35
Set-TextCell $ws 'D35' '5.36'
Set-TextCell $ws 'E35' '  -2.31%  '

# Row 36
Set-TextCell $ws 'D36' '0.411'
Set-TextCell $ws 'E36' '  -3.13%  '

# Row 37
Set-TextCell $ws 'D37' '19.99'
Set-TextCell $ws 'E37' '  -3.02%  '

# Row 38
Set-TextCell $ws 'D38' '0.999'
Set-TextCell $ws 'E38' '  -0.02%  '

# Row 39
Set-TextCell $ws 'D39' '1.91'
Set-TextCell $ws 'E39' '  -1.65%  '

# Row 40
Set-TextCell $ws 'D40' '151.35'
Set-TextCell $ws 'E40' '  -4.34%  '

# Row 41
Set-TextCell $ws 'D41' '0.999'
Set-TextCell $ws 'E41' '  -0.02%  '

# Row 42
Set-TextCell $ws 'B42' 'dogwifhat'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell $ws 'D42' '2.41'
Set-TextCell $ws 'E42' '  +1.05%  '

# Row 43
Set-TextCell $ws 'B43' 'Aave'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D43' '158.53'
Set-TextCell $ws 'E43' '  -2.72%  '

# Row 44
Set-TextCell $ws 'D44' '3.97'
Set-TextCell $ws 'E44' '  -3.71%  '

# Row 45
Set-TextCell $ws 'D45' '0.0602'
Set-TextCell $ws 'E45' '  -1.57%  '

# Row 46
Set-TextCell $ws 'D46' '22.64'
Set-TextCell $ws 'E46' '  -0.30%  '

# Row 47
Set-TextCell $ws 'D47' '0.631'
Set-TextCell $ws 'E47' '  -1.76%  '

# Row 48
Set-TextCell $ws 'D48' '0.101'
Set-TextCell $ws 'E48' '  +2.24%  '

# Row 49
Set-TextCell $ws 'E49' '  -3.20%  '

# Row 50
Set-TextCell $ws 'D50' '19.10'
Set-TextCell $ws 'E50' '  -5.45%  '

# Row 51
Set-TextCell $ws 'E51' '  -7.38%  '
